$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 263.96295
$ws.Range("I33").Value = 213.18182
$ws.Range("K33").Value = 213.18182
$ws.Range("M33").Value = 15.81818000000001

$ws.Range("H97").Value = 338.42856
$ws.Range("J97").Value = 311.5
$ws.Range("L97").Value = 934.5
$ws.Range("N97").Value = -1926.5

$ws.Range("H127").Value = 1312.091
$ws.Range("I127").Value = 619.06665
$ws.Range("J127").Value = 2797.1428
$ws.Range("K127").Value = 1857.19995
$ws.Range("L127").Value = 8391.428400000001
$ws.Range("M127").Value = 3102.80005
$ws.Range("N127").Value = -18311.4284

$ws.Range("H129").Value = 700.625
$ws.Range("I129").Value = 516.375
$ws.Range("J129").Value = 884.875
$ws.Range("K129").Value = 1549.125
$ws.Range("L129").Value = 2654.625
$ws.Range("M129").Value = 3450.875
$ws.Range("N129").Value = -12654.625


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 671.13043
$ws.Range("I2").Value = 453.7857
$ws.Range("J2").Value = 1009.2222
$ws.Range("K2").Value = 453.7857
$ws.Range("L2").Value = 1009.2222
$ws.Range("M2").Value = -340.7857
$ws.Range("N2").Value = -1235.2222

$ws.Range("H32").Value = 4150.125
$ws.Range("I32").Value = 4229.2407
$ws.Range("K32").Value = 4229.2407
$ws.Range("M32").Value = -3942.2407

$ws.Range("H39").Value = 0
$ws.Range("I39").Value = 0
$ws.Range("K39").Value = 0
$ws.Range("M39").ClearContents()

$ws.Range("H41").Value = 6609.3335
$ws.Range("I41").Value = 5131.2
$ws.Range("K41").Value = 5131.2
$ws.Range("M41").Value = -4717.2

$ws.Range("H45").Value = 1190.421
$ws.Range("I45").Value = 1123.7646
$ws.Range("J45").Value = 1757
$ws.Range("K45").Value = 1123.7646
$ws.Range("L45").Value = 1757
$ws.Range("M45").Value = -746.7646
$ws.Range("N45").Value = -2511

$ws.Range("H63").Value = 2066.6592
$ws.Range("I63").Value = 1930.4482
$ws.Range("K63").Value = 1930.4482
$ws.Range("M63").Value = -1244.4482

$ws.Range("H66").Value = 2066.6592
$ws.Range("I66").Value = 1930.4482
$ws.Range("K66").Value = 9652.241
$ws.Range("M66").Value = -6220.241

$ws.Range("H110").Value = 1804.1428
$ws.Range("I110").Value = 578.8333
$ws.Range("J110").Value = 2723.125
$ws.Range("K110").Value = 578.8333
$ws.Range("L110").Value = 2723.125
$ws.Range("M110").Value = 1466.1667
$ws.Range("N110").Value = -6813.125

$ws.Range("H116").Value = 671.13043
$ws.Range("I116").Value = 453.7857
$ws.Range("J116").Value = 1009.2222
$ws.Range("K116").Value = 453.7857
$ws.Range("L116").Value = 1009.2222
$ws.Range("M116").Value = 1840.2143
$ws.Range("N116").Value = -5597.2222

$ws.Range("H122").Value = 1286.5555
$ws.Range("I122").Value = 1286.5555
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 3859.6665
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -1409.6665
$ws.Range("N122").ClearContents()


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 671.13043
$ws.Range("I3").Value = 453.7857
$ws.Range("J3").Value = 1009.2222
$ws.Range("K3").Value = 453.7857
$ws.Range("L3").Value = 1009.2222
$ws.Range("M3").Value = -339.7857
$ws.Range("N3").Value = -1237.2222

$ws.Range("H80").Value = 491.6842
$ws.Range("J80").Value = 566.6875
$ws.Range("L80").Value = 566.6875
$ws.Range("N80").Value = -2562.6875

$ws.Range("H83").Value = 491.6842
$ws.Range("J83").Value = 566.6875
$ws.Range("L83").Value = 2833.4375
$ws.Range("N83").Value = -12817.4375

$ws.Range("H99").Value = 200001140
$ws.Range("J99").Value = 1345
$ws.Range("L99").Value = 1345
$ws.Range("N99").Value = -4341

$ws.Range("H105").Value = 77685060
$ws.Range("I105").Value = 91809370
$ws.Range("J105").Value = 1349
$ws.Range("K105").Value = 91809370
$ws.Range("L105").Value = 1349
$ws.Range("M105").Value = -91807623
$ws.Range("N105").Value = -4843

$ws.Range("H107").Value = 1164.6471
$ws.Range("I107").Value = 888.3333
$ws.Range("K107").Value = 888.3333
$ws.Range("M107").Value = 1031.6667

$ws.Range("H134").Value = 7083.684
$ws.Range("I134").Value = 1211.4
$ws.Range("J134").Value = 13608.444
$ws.Range("K134").Value = 3634.2
$ws.Range("L134").Value = 40825.33199999999
$ws.Range("M134").Value = -1099.2
$ws.Range("N134").Value = -45895.33199999999


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1443.1578
$ws.Range("I31").Value = 1391.3334
$ws.Range("K31").Value = 1391.3334
$ws.Range("M31").Value = -1096.3334

$ws.Range("H34").Value = 1443.1578
$ws.Range("I34").Value = 1391.3334
$ws.Range("K34").Value = 1391.3334
$ws.Range("M34").Value = -1189.3334

$ws.Range("H52").Value = 41713.332
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 41713.332
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 41713.332
$ws.Range("M52").ClearContents()
$ws.Range("N52").Value = -42301.332

$ws.Range("H132").Value = 2127.724
$ws.Range("I132").Value = 1864.8235
$ws.Range("J132").Value = 2500.1667
$ws.Range("K132").Value = 5594.470499999999
$ws.Range("L132").Value = 7500.500100000001
$ws.Range("M132").Value = -3064.470499999999
$ws.Range("N132").Value = -12560.5001

$ws.Range("H134").Value = 33335990
$ws.Range("I134").Value = 3104.4546
$ws.Range("K134").Value = 9313.363799999999
$ws.Range("M134").Value = -6778.363799999999


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 615.8158
$ws.Range("I113").Value = 533.1539
$ws.Range("J113").Value = 658.8
$ws.Range("K113").Value = 1599.4617
$ws.Range("L113").Value = 1976.4
$ws.Range("M113").Value = 570.5382999999999
$ws.Range("N113").Value = -6316.4

$ws.Range("H131").Value = 16952340
$ws.Range("J131").Value = 3633.5881
$ws.Range("L131").Value = 10900.7643
$ws.Range("N131").Value = -20980.7643

$ws.Range("H132").Value = 864.2
$ws.Range("I132").Value = 763.1429000000001
$ws.Range("J132").Value = 1100
$ws.Range("K132").Value = 6868.2861
$ws.Range("L132").Value = 9900
$ws.Range("M132").Value = -4338.2861
$ws.Range("N132").Value = -14960

$ws.Range("H140").Value = 43204.36
$ws.Range("I140").Value = 94291.73
$ws.Range("J140").Value = 3064.2856
$ws.Range("K140").Value = 282875.19
$ws.Range("L140").Value = 9192.856800000001
$ws.Range("M140").Value = -277695.19
$ws.Range("N140").Value = -19552.8568


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2077.8125
$ws.Range("I126").Value = 1775.5454
$ws.Range("J126").Value = 2742.8
$ws.Range("K126").Value = 5326.6362
$ws.Range("L126").Value = 8228.400000000001
$ws.Range("M126").Value = -2856.6362
$ws.Range("N126").Value = -13168.4

$ws.Range("H136").Value = 24285.572
$ws.Range("J136").Value = 24285.572
$ws.Range("L136").Value = 72856.716
$ws.Range("N136").Value = -77956.716


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2828.75
$ws.Range("I40").Value = 1894.25
$ws.Range("J40").Value = 7501.25
$ws.Range("K40").Value = 1894.25
$ws.Range("L40").Value = 7501.25
$ws.Range("M40").Value = -1758.25
$ws.Range("N40").Value = -7773.25

$ws.Range("H93").Value = 975.8570999999999
$ws.Range("I93").Value = 926.2
$ws.Range("K93").Value = 926.2
$ws.Range("M93").Value = 321.8

$ws.Range("H122").Value = 41668324
$ws.Range("I122").Value = 62501076
$ws.Range("J122").Value = 2827.5
$ws.Range("K122").Value = 187503228
$ws.Range("L122").Value = 8482.5
$ws.Range("M122").Value = -187500778
$ws.Range("N122").Value = -13382.5


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 400
$ws.Range("I100").Value = 400
$ws.Range("K100").Value = 800
$ws.Range("M100").Value = -259

$ws.Range("H113").Value = 346.78262
$ws.Range("I113").Value = 266.27777
$ws.Range("J113").Value = 636.6
$ws.Range("K113").Value = 798.83331
$ws.Range("L113").Value = 1909.8
$ws.Range("M113").Value = 1371.16669
$ws.Range("N113").Value = -6249.8

$ws.Range("H136").Value = 934.2222
$ws.Range("I136").Value = 856.8461
$ws.Range("J136").Value = 1135.4
$ws.Range("K136").Value = 2570.5383
$ws.Range("L136").Value = 3406.2
$ws.Range("M136").Value = -20.53830000000016
$ws.Range("N136").Value = -8506.200000000001

